# Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)
# Updates the metrics computed for row 3 (file index 1 / second data row)
# in the metrics_sim.xlsx results table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C3"  = 0.8888888888888888
    "H3"  = 0.7834031214312904
    "I3"  = 0.03709401709401709
    "J3"  = 0.7777777777777778
    "K3"  = 110.4444444444444
    "Q3"  = 3
    "R3"  = 12
    "S3"  = 32
    "T3"  = 92
    "U3"  = 219
    "V3"  = 2597
    "W3"  = 2588
    "X3"  = 2568
    "Y3"  = 2508
    "Z3"  = 2381
    "AF3" = 0.998846
    "AG3" = 0.995385
    "AH3" = 0.987692
    "AI3" = 0.964615
    "AJ3" = 0.915769
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
